# Update the "dSF" (column F) values for the rows listed in the diff.
# Mapping of worksheet row number -> new value for column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    3  = -5
    9  = -1
    12 = 3
    13 = -2
    16 = -3
    20 = 5
    25 = -2
    29 = -1
    30 = -1
    37 = -4
    38 = 2
    45 = -5
    46 = 0
    55 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
